# aoiConditions/train2P2Block4Test.xlsx: "updated condiitons for none"
#
# The condition column (C) is being repurposed: the header is renamed from
# "audioFalse" to "currentPhase", and the per-row audio-file values are
# replaced with the training-phase label "train2P2". The sibling column D
# (image file names) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: C1 "audioFalse" -> "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Data rows: C2/C3 old audio filenames -> "train2P2"
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
